$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell content (keep styles/formatting intact)
$ws.Cells.ClearContents()

# Row 1: column index header (0-18) across B1:T1
$colIdx = 0
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")) {
    $ws.Range("$col" + "1").Value = $colIdx
    $colIdx++
}

# Column A (sample index) and column B (sample name) for rows 2-31
# Writing column B top-to-bottom first builds the shared-string table
# in the exact order required (indices 0-29).
$names = @("HKL", "Spiral5", "Holden", "Rizzie Spiral", "RotRing OmegaMax-90", "Equal Angle", "Tilt Rotate", "CLR", "Rizzie Hex", "Matthies Hex", "Tilt Rotate_Partial", "RotRing OmegaMax-60", "Equal Angle_Partial", "Rizzie Hex_Partial", "ND Single", "RD Single", "TD Single", "Morris Single", "Ring Perpendicular to ND", "Ring Perpendicular to RD", "Ring Perpendicular to TD", "OffsetFTD", "OffsetATD", "OffsetF45", "OffsetA45", "OffsetFRD", "OffsetARD", "Gaussian Quadrature", "Michael-CCHex", "Michael-SNHex")
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
}

# Rows 30-31 are brand new (beyond the original A1:T29 extent), so they
# start with no cell formatting at all. Clone the "index" column style
# (bold/centered/bordered) from row 29 so A30:A31 pick up the same xf
# as every other row in column A, without minting a new cellXfs entry.
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 2 (header) labels across C2:T2 -- appends shared-string
# indices 30-47 in left-to-right order.
$labels = @("[3, 2, 1]", "[1, 1, 0]", "[2, 2, 2]", "[3, 1, 0]", "[2, 2, 0]", "[2, 0, 0]", "[2, 1, 1]", "[4, 0, 0]", "1Pair-A", "1Pair-B", "2Pairs-A", "2Pairs-B", "3Pairs-A", "3Pairs-B", "3Pairs-C", "4Pairs", "5A4F", "MaxUnique")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(2, $i + 3).Value = $labels[$i]
}

# Data body: C3:T31 filled with 1
for ($r = 3; $r -le 31; $r++) {
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
